$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new values
$ws.Range("A2").Value = 45034.50694444445
$ws.Range("B2").Value = 7.205
$ws.Range("C2").Value = 5.605
$ws.Range("D2").Value = 1.192
$ws.Range("E2").Value = 15.61
$ws.Range("F2").Value = 12.374
$ws.Range("G2").Value = 4.794
$ws.Range("H2").Value = 14.798
$ws.Range("I2").Value = 8.957000000000001
$ws.Range("J2").Value = 4.429
$ws.Range("K2").Value = 5.631
$ws.Range("L2").Value = 6.249
$ws.Range("M2").Value = 7.306
$ws.Range("N2").Value = 2.788
$ws.Range("O2").Value = 6.015
$ws.Range("P2").Value = 7.738
$ws.Range("Q2").Value = 5.138
$ws.Range("R2").Value = 0.492
$ws.Range("S2").Value = 0.931
$ws.Range("T2").Value = 84.14
$ws.Range("U2").Value = 16.424
$ws.Range("V2").Value = 5.552
$ws.Range("W2").Value = 10.174
$ws.Range("X2").Value = 6.283
$ws.Range("Y2").Value = 0.894
$ws.Range("Z2").Value = 9.741
$ws.Range("AA2").Value = 4.361
$ws.Range("AB2").Value = 4.885
$ws.Range("AC2").Value = 6.06
$ws.Range("AD2").Value = 8.134
$ws.Range("AE2").Value = 1.522
$ws.Range("AF2").Value = 13.18
$ws.Range("AG2").Value = 3.625
$ws.Range("AH2").Value = 6.387

$ws.Range("A3").Value = 45034.51388888889
$ws.Range("B3").Value = 5.966
$ws.Range("C3").Value = 4.535
$ws.Range("D3").Value = 0.599
$ws.Range("E3").Value = 13.216
$ws.Range("F3").Value = 10.406
$ws.Range("G3").Value = 4.154
$ws.Range("H3").Value = 17.694
$ws.Range("I3").Value = 7.333
$ws.Range("J3").Value = 3.732
$ws.Range("K3").Value = 4.605
$ws.Range("L3").Value = 5.247
$ws.Range("M3").Value = 5.938
$ws.Range("N3").Value = 2.095
$ws.Range("O3").Value = 4.888
$ws.Range("P3").Value = 6.522
$ws.Range("Q3").Value = 4.252
$ws.Range("R3").Value = 0.199
$ws.Range("S3").Value = 0.532
$ws.Range("T3").Value = 67.062
$ws.Range("U3").Value = 13.698
$ws.Range("V3").Value = 4.511
$ws.Range("W3").Value = 8.756
$ws.Range("X3").Value = 5.087
$ws.Range("Y3").Value = 0.671
$ws.Range("Z3").Value = 10.066
$ws.Range("AA3").Value = 3.686
$ws.Range("AB3").Value = 3.9
$ws.Range("AC3").Value = 4.666
$ws.Range("AD3").Value = 6.215
$ws.Range("AE3").Value = 0.707
$ws.Range("AF3").Value = 16.639
$ws.Range("AG3").Value = 2.848
$ws.Range("AH3").Value = 5.338

$ws.Range("A4").Value = 45034.52083333334
$ws.Range("B4").Value = 19.981
$ws.Range("C4").Value = 15.027
$ws.Range("D4").Value = 0.921
$ws.Range("E4").Value = 43.703
$ws.Range("F4").Value = 35.711
$ws.Range("G4").Value = 15.331
$ws.Range("H4").Value = 57.172
$ws.Range("I4").Value = 24.262
$ws.Range("J4").Value = 11.249
$ws.Range("K4").Value = 16.01
$ws.Range("L4").Value = 17.493
$ws.Range("M4").Value = 18.773
$ws.Range("N4").Value = 5.447
$ws.Range("O4").Value = 15.79
$ws.Range("P4").Value = 22.226
$ws.Range("Q4").Value = 13.256
$ws.Range("R4").Value = 0.244
$ws.Range("S4").Value = 0.825
$ws.Range("T4").Value = 233.017
$ws.Range("U4").Value = 44.069
$ws.Range("V4").Value = 14.575
$ws.Range("W4").Value = 29.486
$ws.Range("X4").Value = 15.857
$ws.Range("Y4").Value = 2.084
$ws.Range("Z4").Value = 29.373
$ws.Range("AA4").Value = 12.668
$ws.Range("AB4").Value = 11.601
$ws.Range("AC4").Value = 13.678
$ws.Range("AD4").Value = 18.82
$ws.Range("AE4").Value = 0.461
$ws.Range("AF4").Value = 52.188
$ws.Range("AG4").Value = 8.507999999999999
$ws.Range("AH4").Value = 18.015

$ws.Range("A5").Value = 45034.52777777778
$ws.Range("B5").Value = 21.47
$ws.Range("C5").Value = 16.13
$ws.Range("D5").Value = 0.9
$ws.Range("E5").Value = 46.91
$ws.Range("F5").Value = 38.44
$ws.Range("G5").Value = 16.58
$ws.Range("H5").Value = 65.77
$ws.Range("I5").Value = 26.04
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 17.24
$ws.Range("L5").Value = 18.79
$ws.Range("M5").Value = 20.08
$ws.Range("N5").Value = 5.73
$ws.Range("O5").Value = 16.92
$ws.Range("P5").Value = 23.93
$ws.Range("Q5").Value = 14.17
$ws.Range("R5").Value = 0.23
$ws.Range("S5").Value = 0.8100000000000001
$ws.Range("T5").Value = 250.2
$ws.Range("U5").Value = 47.32
$ws.Range("V5").Value = 15.62
$ws.Range("W5").Value = 31.77
$ws.Range("X5").Value = 16.95
$ws.Range("Y5").Value = 2.23
$ws.Range("Z5").Value = 32.73
$ws.Range("AA5").Value = 13.64
$ws.Range("AB5").Value = 12.35
$ws.Range("AC5").Value = 14.54
$ws.Range("AD5").Value = 20.06
$ws.Range("AE5").Value = 0.34
$ws.Range("AF5").Value = 60.05
$ws.Range("AG5").Value = 9.06
$ws.Range("AH5").Value = 19.37

# Delete row 6 (dataset now has 4 data rows instead of 5)
$ws.Rows.Item(6).Delete()

# Adjust column widths that changed from 7 to 8 characters
$ws.Columns.Item(2).ColumnWidth = 7.142857
$ws.Columns.Item(3).ColumnWidth = 7.142857
$ws.Columns.Item(7).ColumnWidth = 7.142857
$ws.Columns.Item(10).ColumnWidth = 7.142857
$ws.Columns.Item(12).ColumnWidth = 7.142857
$ws.Columns.Item(13).ColumnWidth = 7.142857
$ws.Columns.Item(16).ColumnWidth = 7.142857
$ws.Columns.Item(17).ColumnWidth = 7.142857
$ws.Columns.Item(22).ColumnWidth = 7.142857
$ws.Columns.Item(24).ColumnWidth = 7.142857
$ws.Columns.Item(27).ColumnWidth = 7.142857
$ws.Columns.Item(28).ColumnWidth = 7.142857
$ws.Columns.Item(29).ColumnWidth = 7.142857
$ws.Columns.Item(34).ColumnWidth = 7.142857
